$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Swap the full data (columns B..AC) between the following row pairs.
#    These pairs represent two match rows whose contents were transposed
#    (the "id" in column A and the row position itself stay put).
#    A scratch row (far below the used range) is used as a temporary
#    holding area for the swap.
# ---------------------------------------------------------------------------

$pairs = @(
    @(9,10),
    @(31,32),
    @(63,64),
    @(70,71),
    @(84,85),
    @(87,88),
    @(125,126),
    @(142,143),
    @(159,160),
    @(166,167),
    @(183,184),
    @(198,199)
)

$scratchRow = 5000

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $ws.Range("B$r1`:AC$r1").Copy() | Out-Null
    $ws.Range("B$scratchRow`:AC$scratchRow").PasteSpecial(-4163) | Out-Null

    $ws.Range("B$r2`:AC$r2").Copy() | Out-Null
    $ws.Range("B$r1`:AC$r1").PasteSpecial(-4163) | Out-Null

    $ws.Range("B$scratchRow`:AC$scratchRow").Copy() | Out-Null
    $ws.Range("B$r2`:AC$r2").PasteSpecial(-4163) | Out-Null

    $ws.Range("B$scratchRow`:AC$scratchRow").Clear() | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Append two new fixture rows (203 and 204) at the end of the sheet,
#    copying the number/cell formatting from the last existing row (202)
#    for the "id" (A) and "Date" (E) columns.
# ---------------------------------------------------------------------------

$ws.Range("A202").Copy() | Out-Null
$ws.Range("A203").PasteSpecial(-4122) | Out-Null
$ws.Range("A203").PasteSpecial(-4122) | Out-Null

$ws.Range("E202").Copy() | Out-Null
$ws.Range("E203:E204").PasteSpecial(-4122) | Out-Null

$ws.Range("A202:A202").Copy() | Out-Null
$ws.Range("A204").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

$ws.Range("A203").Value = 201
$ws.Range("B203").Value = 6798879
$ws.Range("C203").Value = "Denmark Division 1"
$ws.Range("D203").Value = "Denmark Division 1"
$ws.Range("E203").Value = Get-Date -Year 2024 -Month 2 -Day 24 -Hour 9 -Minute 0 -Second 0
$ws.Range("F203").Value = "AaB"
$ws.Range("G203").Value = "Sonderjyske"
$ws.Range("K203").Value = 2.25
$ws.Range("L203").Value = 3.4
$ws.Range("M203").Value = 2.9
$ws.Range("N203").Value = 2.05
$ws.Range("O203").Value = 3.4
$ws.Range("P203").Value = 3.3
$ws.Range("Q203").Value = -0.25
$ws.Range("R203").Value = 1.8
$ws.Range("S203").Value = 2.05
$ws.Range("T203").Value = 2.75
$ws.Range("U203").Value = 1.85
$ws.Range("V203").Value = 2
$ws.Range("W203").Value = 0
$ws.Range("X203").Value = 0
$ws.Range("Y203").Value = 0
$ws.Range("Z203").Value = 0
$ws.Range("AA203").Value = 0

$ws.Range("A204").Value = 202
$ws.Range("B204").Value = 6800816
$ws.Range("C204").Value = "Denmark Division 1"
$ws.Range("D204").Value = "Denmark Division 1"
$ws.Range("E204").Value = Get-Date -Year 2024 -Month 2 -Day 24 -Hour 11 -Minute 0 -Second 0
$ws.Range("F204").Value = "B93 Copenhagen"
$ws.Range("G204").Value = "Naestved"
$ws.Range("K204").Value = 2.25
$ws.Range("L204").Value = 3.5
$ws.Range("M204").Value = 2.9
$ws.Range("N204").Value = 2.375
$ws.Range("O204").Value = 3.5
$ws.Range("P204").Value = 2.7
$ws.Range("Q204").Value = 0
$ws.Range("R204").Value = 1.8
$ws.Range("S204").Value = 2.05
$ws.Range("T204").Value = 2.75
$ws.Range("U204").Value = 1.975
$ws.Range("V204").Value = 1.875
$ws.Range("W204").Value = 0
$ws.Range("X204").Value = 0
$ws.Range("Y204").Value = 0
$ws.Range("Z204").Value = 0
$ws.Range("AA204").Value = 0
